# Commit: "Start redesign of ADC board"
# Adds a new "Compatible parts can be used?" Y/N column between the
# existing "Note" column (H->I) data and shifts the old Note column
# (I) one slot to the right (J).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1. Insert a new column at I - this pushes the existing "Note"
#    column (and its data/format) from I to J automatically.
# ---------------------------------------------------------------
$ws.Columns.Item(9).Insert()

# Give the freshly inserted column roughly the same auto-sized width
# the original workbook ends up with (best achievable through the
# character-width-based ColumnWidth property).
$ws.Columns.Item(9).ColumnWidth = 33.86

# ---------------------------------------------------------------
# 2. Header + data for the new "Compatible parts can be used?" column.
# ---------------------------------------------------------------
$ws.Range("I1").Value = "Compatible parts can be used?"

$ws.Range("I2").Value = "Y"
$ws.Range("I3").Value = "Y"
$ws.Range("I12").Value = "Y"
$ws.Range("I21").Value = "Y"
$ws.Range("I22").Value = "Y"

# Rows that need the bold "header-like" font (fontId 19: bold, family 3)
# applied along with the border already inherited from the column.
$ws.Range("I4").Value = "N"
$ws.Range("I4").Font.Family = 3
$ws.Range("I4").Font.Bold = $true

# Copy that exact cell format onto the remaining bold N/Y cells so the
# engine reuses the same cell style instead of minting new ones.
$ws.Range("I4").Copy()
$ws.Range("I5:I11").PasteSpecial(-4122)
$ws.Range("I13:I20").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I5").Value = "N"
$ws.Range("I6").Value = "N"
$ws.Range("I7").Value = "N"
$ws.Range("I8").Value = "N"
$ws.Range("I9").Value = "N"
$ws.Range("I10").Value = "N"
$ws.Range("I11").Value = "N"
$ws.Range("I13").Value = "N"
$ws.Range("I14").Value = "N"
$ws.Range("I15").Value = "N"
$ws.Range("I16").Value = "N"
$ws.Range("I17").Value = "N"
$ws.Range("I18").Value = "N"
$ws.Range("I19").Value = "N"
$ws.Range("I20").Value = "N"

# ---------------------------------------------------------------
# 3. Restore the active cell/selection like the authored workbook.
# ---------------------------------------------------------------
$ws.Range("I14").Select()
